$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 286: fill in "Custos do dia" (O286) ---
$ws.Range("O286").Value = 20

# --- Row 290: correct sale value (J290) ---
$ws.Range("J290").Value = 180

# --- Rows 293-296: new sales entries (previously blank rows) ---
$ws.Range("A293").Value = 45787
$ws.Range("B293").Value = "Loja2"
$ws.Range("C293").Value = 2188
$ws.Range("D293").Value = 39
$ws.Range("E293").Value = "Netony"
$ws.Range("F293").Value = "Caterpilhar Cano Curto"
$ws.Range("G293").Value = 100
$ws.Range("H293").Value = 199
$ws.Range("I293").Value = 1
$ws.Range("J293").Value = 190

$ws.Range("A294").Value = 45787
$ws.Range("B294").Value = "Loja2"
$ws.Range("C294").Value = 1046
$ws.Range("D294").Value = 39
$ws.Range("E294").Value = "Rossanfort"
$ws.Range("F294").Value = "Sapatilha"
$ws.Range("G294").Value = 50
$ws.Range("H294").Value = 110
$ws.Range("I294").Value = 1
$ws.Range("J294").Value = 100

$ws.Range("A295").Value = 45787
$ws.Range("B295").Value = "Loja2"
$ws.Range("C295").Value = 1046
$ws.Range("D295").Value = 40
$ws.Range("E295").Value = "Rossanfort"
$ws.Range("F295").Value = "Sapatilha"
$ws.Range("G295").Value = 50
$ws.Range("H295").Value = 110
$ws.Range("I295").Value = 1
$ws.Range("J295").Value = 100

$ws.Range("A296").Value = 45787
$ws.Range("B296").Value = "Loja2"
$ws.Range("C296").Value = 1046
$ws.Range("D296").Value = 43
$ws.Range("E296").Value = "Rossanfort"
$ws.Range("F296").Value = "Sapatilha"
$ws.Range("G296").Value = 50
$ws.Range("H296").Value = 110
$ws.Range("I296").Value = 1
$ws.Range("J296").Value = 100

# --- Row 306: new "Credito"/pagamento tracking entry ---
$ws.Range("W306").Value = 160
$ws.Range("X306").Value = "3x tênis "

# --- Update view: scroll/select near the new data ---
$ws.Range("X306").Select()
